$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI Dll1-Notch4 results (ECs now included as a sending cluster,
# per Dr Hou advice) - replaces the previous 3-row table with 6 rows.
$data = @(
    @("ECs", "Dll1", "Notch4", "ECs", 2, 0.6666666666666666, 10.34874566666667, 31.046237, 0.614862320492409, 0.6148623204924089, 3, 1, 29.546424, 88.639272, 0.9033225104610835, 0.9033225104610834, 305.768427335496, 2751.915846019464, 0.5554189749351301, 0.55541897493513),
    @("ECs", "Dll1", "Notch4", "FAPs", 2, 0.6666666666666666, 10.34874566666667, 31.046237, 0.614862320492409, 0.6148623204924089, 3, 1, 1.672785333333334, 5.018356000000001, 0.05114204841740398, 0.05114204841740398, 17.31122996959689, 155.801069726372, 0.03144531856466014, 0.03144531856466014),
    @("ECs", "Dll1", "Notch4", "sCs", 2, 0.6666666666666666, 10.34874566666667, 31.046237, 0.614862320492409, 0.6148623204924089, 3, 1, 1.489401, 4.468203, 0.04553544112151264, 0.04553544112151264, 15.413432144679, 138.720889302111, 0.02799802699261873, 0.02799802699261872),
    @("sCs", "Dll1", "Notch4", "ECs", 3, 1, 6.482251000000001, 19.446753, 0.3851376795075911, 0.385137679507591, 3, 1, 29.546424, 88.639272, 0.9033225104610835, 0.9033225104610834, 191.527336520424, 1723.746028683816, 0.3479035355259534, 0.3479035355259533),
    @("sCs", "Dll1", "Notch4", "FAPs", 3, 1, 6.482251000000001, 19.446753, 0.3851376795075911, 0.385137679507591, 3, 1, 1.672785333333334, 5.018356000000001, 0.05114204841740398, 0.05114204841740398, 10.84341439978534, 97.59072959806802, 0.01969672985274384, 0.01969672985274384),
    @("sCs", "Dll1", "Notch4", "sCs", 3, 1, 6.482251000000001, 19.446753, 0.3851376795075911, 0.385137679507591, 3, 1, 1.489401, 4.468203, 0.04553544112151264, 0.04553544112151264, 9.654671121651, 86.892040094859, 0.01753741412889392, 0.01753741412889391),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $row[$j]
    }
}
